$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'317.73"
$ws.Range("E2").Value = "'3.89%"
$ws.Range("D3").Value = "'39.66"
$ws.Range("E3").Value = "'1.89%"
$ws.Range("E4").Value = "'0.62%"
$ws.Range("D5").Value = "'0.08205"
$ws.Range("D6").Value = "'2.040"
$ws.Range("E6").Value = "'5.97%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.342"
$ws.Range("E7").Value = "'3.45%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'8.338"
$ws.Range("E8").Value = "'4.21%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9382"
$ws.Range("E9").Value = "'1.06%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1354"
$ws.Range("E10").Value = "'-8.88%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1990"
$ws.Range("E11").Value = "'3.69%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09120"
$ws.Range("E12").Value = "'1.23%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03498"
$ws.Range("E13").Value = "'-0.62%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09792"
$ws.Range("E14").Value = "'0.20%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001409"
$ws.Range("E15").Value = "'1.25%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006192"
$ws.Range("E16").Value = "'5.31%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.693"
$ws.Range("E17").Value = "'-2.33%"
$ws.Range("D18").Value = "'3.290"
$ws.Range("E18").Value = "'-3.55%"
$ws.Range("E19").Value = "'1.39%"
$ws.Range("E20").Value = "'-0.51%"
$ws.Range("D21").Value = "'4.957"
$ws.Range("E21").Value = "'5.94%"
$ws.Range("D22").Value = "'0.2452"
$ws.Range("E22").Value = "'1.43%"
$ws.Range("D23").Value = "'0.04365"
$ws.Range("E23").Value = "'-0.45%"
$ws.Range("E24").Value = "'-0.40%"
$ws.Range("D25").Value = "'0.004794"
$ws.Range("E25").Value = "'12.22%"
$ws.Range("E26").Value = "'-0.07%"
$ws.Range("D27").Value = "'0.0004007"
$ws.Range("E27").Value = "'-9.91%"
$ws.Range("D39").Value = "'0.02249"
$ws.Range("E39").Value = "'10.72%"
$ws.Range("D40").Value = "'0.05187"
$ws.Range("E40").Value = "'2.63%"
$ws.Range("D41").Value = "'0.007755"
$ws.Range("E41").Value = "'3.05%"
$ws.Range("D42").Value = "'0.009959"
$ws.Range("E42").Value = "'2.33%"
$ws.Range("E43").Value = "'4.15%"
$ws.Range("D44").Value = "'0.002051"
$ws.Range("E44").Value = "'-2.54%"
$ws.Range("D45").Value = "'0.009123"
$ws.Range("E45").Value = "'-7.89%"
$ws.Range("D46").Value = "'0.00006612"
$ws.Range("E46").Value = "'6.51%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.04%"
$ws.Range("D48").Value = "'0.002950"
$ws.Range("E48").Value = "'2.69%"
$ws.Range("D49").Value = "'0.001693"
$ws.Range("E49").Value = "'-6.15%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'-0.04%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'-0.04%"
